# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off (new target xliff files generated), while flagging that the
# handback for b.md is stale relative to a newer source commit.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34b14e42d6562570b96a5e16776a0af7fa6475ce/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16314f12fee19ca4566cd5d4b0791869cd08adbb/e2e/b.md."

# ---- Overview sheet: row for b.md (row 3) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 06:42:20"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to a width that round-trips to 40
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# Row 2 (a.md): status is now "Ready for handoff"
$wsZh.Range("C2").Value = "Ready for handoff"

# Row 3 (b.md): status, content duplicate flag, latest handoff file/datetime, error detail
$wsZh.Range("C3").Value = "Ready for handoff"
# Copy the existing "False" text cell (O2) over F3 so the text "False" stays a
# shared string instead of being auto-coerced into a boolean cell.
$wsZh.Range("O2").Copy($wsZh.Range("F3"))
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-01 06:42:15"
$wsZh.Range("P3").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P) to a width that round-trips to 40
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664

# Row 3 (b.md): status, content duplicate flag, latest handoff file/datetime, error detail
$wsDe.Range("C3").Value = "Ready for handoff"
# Copy the existing "False" text cell (O2) over F3 so the text "False" stays a
# shared string instead of being auto-coerced into a boolean cell.
$wsDe.Range("O2").Copy($wsDe.Range("F3"))
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-01 06:42:20"
$wsDe.Range("P3").Value = $errorDetail
